$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Year of Treatment" column (column B). All subsequent
# columns (All substances, Opioids, Cocaine, Stimulants,
# Hypnotics and Sedatives, Hallucinogens, Volatile Inhalants, Cannabis,
# Other substances) shift one position to the left, from C:K to B:J.
$ws.Range("B:B").Delete()
